$wb = $excel.ActiveWorkbook

# Week 17 roster update: add newly-signed player L.Bell as a new tracked
# player column (inserted right after K.Vaughn) on both the Rushing and
# Receiving yards-list sheets, defaulting his "logged" flag to "n" like
# every other player.
foreach ($ws in $wb.Worksheets) {
    $ws.Columns("I:I").Insert()
    $ws.Range("I1").Value = "L.Bell"
    $ws.Range("I2").Value = "n"
}
